# feat: Enable AgGrid column grouping, and add import/export in customToolbar
#
# Data-level changes to cypress/fixtures/customers.xlsx:
#  - "customers" sheet (sheet1) becomes/stays the active sheet/tab
#  - Howard's birthday (E5) changes from 12/05/1987 to 21/05/2002
#  - A new row is appended for Billy: age 23, canDrinkAlcohol TRUE,
#    favoriteDrink Beer, birthday DATE(1940,4,28), height 1.25

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("customers")

# Update Howard's birthday
$ws.Range("E5").Formula = "21/05/2002"

# Append new row 6 for Billy
$ws.Range("A6").Value = "Billy"
$ws.Range("B6").Value = 23
$ws.Range("C6").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = "Beer"
$ws.Range("E6").NumberFormat = "mm/dd/yy"
$ws.Range("E6").Formula = "=DATE(1940,4,28)"
$ws.Range("F6").Value = 1.25

# Make the customers sheet the active/selected one, with E7 selected
$ws.Activate() | Out-Null
$ws.Range("E7").Select() | Out-Null
